# "added width and height spec"
#
# On the "ADC_100MS" sheet, insert two new rows of data (width / height)
# above the existing "pin" table, pushing the whole table (and its blank
# spacer rows) down by three rows. Row 1 (A1) is left untouched, the new
# data lands on rows 3-4, and a blank row 2 remains between them (mirroring
# the existing blank-spacer-row convention already used in this table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADC_100MS")

# Insert 3 blank rows starting at row 2 - this shifts the old row 2 (blank),
# row 3..13 data/formulas down to row 5..16, preserving relative formulas.
$ws.Rows.Item(2).Resize(3).Insert()

# The inserted rows inherit row 1's formatting in column A; clear that so no
# stray formatted-but-empty cells are left behind in column A.
$ws.Range("A2:A4").Clear()

# Fill in the new width/height spec rows.
$ws.Range("B3").Value = "width"
$ws.Range("C3").Value = 500
$ws.Range("B4").Value = "height"
$ws.Range("C4").Value = 300

# Match the author's resulting selection.
$ws.Range("C5").Select()
